# Daily_Updates.xlsx — 2025-09-01 OneDrive auto-sync edit
#
# The only substantive data change in this commit is on the "Daily Updates"
# sheet: the values that had been entered in D2 (CERTIFICATIONS row's
# "CDF" note) and I2 (the DUE_DATE 45896) were cleared out again, leaving
# the rest of row 2 (E2:H2 = HIGH / IN PROCESS / GK / NKC) untouched.
# Clearing D2 also drops "CDF" from the shared-string table on save, which
# cascades the shared-string indices used by every other sheet — that is
# an automatic side effect of the save, not something to replicate by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily Updates")

$ws.Range("D2").ClearContents()
$ws.Range("I2").ClearContents()

# Reflect the cursor ending up on I2 (matches the saved selection in the
# sheetView) after clearing the row.
$ws.Activate()
$ws.Range("I2").Select() | Out-Null
